$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I15").Value = "sd"
$ws.Range("J15").Value = "Statement-non-opinion"
$ws.Range("I21").Value = "sd"
$ws.Range("J21").Value = "Statement-non-opinion"
$ws.Range("I29").Value = "sd"
$ws.Range("J29").Value = "Statement-non-opinion"
$ws.Range("I38").Value = "sd"
$ws.Range("J38").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "b"
$ws.Range("J39").Value = "Acknowledge (Backchannel)"
$ws.Range("I45").Value = "b"
$ws.Range("J45").Value = "Acknowledge (Backchannel)"
$ws.Range("I54").Value = "aa"
$ws.Range("J54").Value = "Agree/Accept"
$ws.Range("I55").Value = "aa"
$ws.Range("J55").Value = "Agree/Accept"
$ws.Range("I67").Value = "ba"
$ws.Range("J67").Value = "Appreciation"
$ws.Range("I99").Value = "sd"
$ws.Range("J99").Value = "Statement-non-opinion"
$ws.Range("I109").Value = "sd"
$ws.Range("J109").Value = "Statement-non-opinion"
$ws.Range("I114").Value = "aa"
$ws.Range("J114").Value = "Agree/Accept"
$ws.Range("I133").Value = "sv"
$ws.Range("J133").Value = "Statement-opinion"
$ws.Range("I136").Value = "ba"
$ws.Range("J136").Value = "Appreciation"
$ws.Range("I140").Value = "ba"
$ws.Range("J140").Value = "Appreciation"
$ws.Range("I145").Value = "%"
$ws.Range("J145").Value = "Uninterpretable"
$ws.Range("I175").Value = "sd"
$ws.Range("J175").Value = "Statement-non-opinion"
$ws.Range("I178").Value = "sd"
$ws.Range("J178").Value = "Statement-non-opinion"
$ws.Range("I179").Value = "aa"
$ws.Range("J179").Value = "Agree/Accept"
$ws.Range("I180").Value = "sd"
$ws.Range("J180").Value = "Statement-non-opinion"
$ws.Range("I184").Value = "sd"
$ws.Range("J184").Value = "Statement-non-opinion"
$ws.Range("I188").Value = "sd"
$ws.Range("J188").Value = "Statement-non-opinion"
$ws.Range("I190").Value = "aa"
$ws.Range("J190").Value = "Agree/Accept"
$ws.Range("I194").Value = "sd"
$ws.Range("J194").Value = "Statement-non-opinion"
$ws.Range("I195").Value = "sv"
$ws.Range("J195").Value = "Statement-opinion"
$ws.Range("I229").Value = "%"
$ws.Range("J229").Value = "Uninterpretable"
$ws.Range("I233").Value = "sd"
$ws.Range("J233").Value = "Statement-non-opinion"
$ws.Range("I242").Value = "aa"
$ws.Range("J242").Value = "Agree/Accept"
$ws.Range("I245").Value = "%"
$ws.Range("J245").Value = "Uninterpretable"
$ws.Range("I246").Value = "aa"
$ws.Range("J246").Value = "Agree/Accept"
$ws.Range("I251").Value = "sv"
$ws.Range("J251").Value = "Statement-opinion"
$ws.Range("I252").Value = "aa"
$ws.Range("J252").Value = "Agree/Accept"
$ws.Range("I263").Value = "aa"
$ws.Range("J263").Value = "Agree/Accept"
$ws.Range("I267").Value = "aa"
$ws.Range("J267").Value = "Agree/Accept"
$ws.Range("I278").Value = "aa"
$ws.Range("J278").Value = "Agree/Accept"
$ws.Range("I284").Value = "aa"
$ws.Range("J284").Value = "Agree/Accept"
$ws.Range("I292").Value = "sd"
$ws.Range("J292").Value = "Statement-non-opinion"
$ws.Range("I293").Value = "sd"
$ws.Range("J293").Value = "Statement-non-opinion"
$ws.Range("I296").Value = "aa"
$ws.Range("J296").Value = "Agree/Accept"
$ws.Range("I298").Value = "sd"
$ws.Range("J298").Value = "Statement-non-opinion"
$ws.Range("I300").Value = "sd"
$ws.Range("J300").Value = "Statement-non-opinion"
$ws.Range("I301").Value = "sd"
$ws.Range("J301").Value = "Statement-non-opinion"
$ws.Range("I310").Value = "aa"
$ws.Range("J310").Value = "Agree/Accept"
$ws.Range("I313").Value = "sd"
$ws.Range("J313").Value = "Statement-non-opinion"
$ws.Range("I314").Value = "aa"
$ws.Range("J314").Value = "Agree/Accept"
$ws.Range("I315").Value = "ba"
$ws.Range("J315").Value = "Appreciation"
$ws.Range("I320").Value = "sd"
$ws.Range("J320").Value = "Statement-non-opinion"
$ws.Range("I339").Value = "ba"
$ws.Range("J339").Value = "Appreciation"
$ws.Range("I345").Value = "sd"
$ws.Range("J345").Value = "Statement-non-opinion"
$ws.Range("I348").Value = "sd"
$ws.Range("J348").Value = "Statement-non-opinion"
$ws.Range("I357").Value = "qy"
$ws.Range("J357").Value = "Yes-No-Question"
$ws.Range("I361").Value = "sd"
$ws.Range("J361").Value = "Statement-non-opinion"
$ws.Range("I368").Value = "sd"
$ws.Range("J368").Value = "Statement-non-opinion"
$ws.Range("I385").Value = "b"
$ws.Range("J385").Value = "Acknowledge (Backchannel)"
$ws.Range("I393").Value = "sd"
$ws.Range("J393").Value = "Statement-non-opinion"
$ws.Range("I424").Value = "sd"
$ws.Range("J424").Value = "Statement-non-opinion"
$ws.Range("I429").Value = "%"
$ws.Range("J429").Value = "Uninterpretable"
$ws.Range("I447").Value = "b"
$ws.Range("J447").Value = "Acknowledge (Backchannel)"
$ws.Range("I456").Value = "sd"
$ws.Range("J456").Value = "Statement-non-opinion"
$ws.Range("I460").Value = "b"
$ws.Range("J460").Value = "Acknowledge (Backchannel)"
$ws.Range("I466").Value = "b"
$ws.Range("J466").Value = "Acknowledge (Backchannel)"
